$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5142.857
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 5142.857
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 5142.857
$ws.Range("N74").Value = -7014.857

$ws.Range("H77").Value = 5142.857
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 5142.857
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 25714.285
$ws.Range("N77").Value = -35074.285

$ws.Range("H98").Value = 1155.1111
$ws.Range("I98").Value = 913.8570999999999
$ws.Range("J98").Value = 1999.5
$ws.Range("K98").Value = 913.8570999999999
$ws.Range("L98").Value = 1999.5
$ws.Range("M98").Value = 584.1429000000001
$ws.Range("N98").Value = -4995.5

$ws.Range("H122").Value = 1155.1111
$ws.Range("I122").Value = 913.8570999999999
$ws.Range("J122").Value = 1999.5
$ws.Range("K122").Value = 2741.5713
$ws.Range("L122").Value = 5998.5
$ws.Range("M122").Value = -291.5712999999996
$ws.Range("N122").Value = -10898.5

$ws.Range("H129").Value = 970.24194
$ws.Range("I129").Value = 365.66666
$ws.Range("J129").Value = 1000.98303
$ws.Range("K129").Value = 1096.99998
$ws.Range("L129").Value = 3002.94909
$ws.Range("M129").Value = 3903.00002
$ws.Range("N129").Value = -13002.94909

$ws.Range("H137").Value = 1254.3429
$ws.Range("I137").Value = 1214.5518
$ws.Range("J137").Value = 1446.6666
$ws.Range("K137").Value = 3643.6554
$ws.Range("L137").Value = 4339.9998
$ws.Range("M137").Value = -1093.6554
$ws.Range("N137").Value = -9439.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4630.4
$ws.Range("I32").Value = 3909.7341
$ws.Range("J32").Value = 9806.091
$ws.Range("K32").Value = 3909.7341
$ws.Range("L32").Value = 9806.091
$ws.Range("M32").Value = -3622.7341
$ws.Range("N32").Value = -10380.091

$ws.Range("H97").Value = 47620270
$ws.Range("I97").Value = 772.2857
$ws.Range("J97").Value = 142859260
$ws.Range("K97").Value = 772.2857
$ws.Range("L97").Value = 142859260
$ws.Range("M97").Value = -276.2857
$ws.Range("N97").Value = -142860252

$ws.Range("H110").Value = 734.4167
$ws.Range("I110").Value = 641.3
$ws.Range("J110").Value = 1200
$ws.Range("K110").Value = 641.3
$ws.Range("L110").Value = 1200
$ws.Range("M110").Value = 1403.7
$ws.Range("N110").Value = -5290

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").ClearContents()

$ws.Range("H131").Value = 60070.6
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 60070.6
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 60070.6
$ws.Range("N131").Value = -70150.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H119").Value = 43770
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 43770
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 43770
$ws.Range("N119").Value = -53446

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1222.9166
$ws.Range("I16").Value = 1051.3334
$ws.Range("J16").Value = 1737.6666
$ws.Range("K16").Value = 1051.3334
$ws.Range("L16").Value = 1737.6666
$ws.Range("M16").Value = -764.3334
$ws.Range("N16").Value = -2311.6666

$ws.Range("H86").Value = 33369346
$ws.Range("I86").Value = 21000
$ws.Range("J86").Value = 55601576
$ws.Range("K86").Value = 21000
$ws.Range("L86").Value = 55601576
$ws.Range("M86").Value = -19877
$ws.Range("N86").Value = -55603822

$ws.Range("H89").Value = 33369346
$ws.Range("I89").Value = 21000
$ws.Range("J89").Value = 55601576
$ws.Range("K89").Value = 105000
$ws.Range("L89").Value = 278007880
$ws.Range("M89").Value = -99384
$ws.Range("N89").Value = -278019112

$ws.Range("H99").Value = 33336580
$ws.Range("I99").Value = 2244.5715
$ws.Range("J99").Value = 62504124
$ws.Range("K99").Value = 2244.5715
$ws.Range("L99").Value = 62504124
$ws.Range("M99").Value = -746.5715
$ws.Range("N99").Value = -62507120

$ws.Range("H105").Value = 5000966
$ws.Range("I105").Value = 6250835
$ws.Range("J105").Value = 1490.2
$ws.Range("K105").Value = 6250835
$ws.Range("L105").Value = 1490.2
$ws.Range("M105").Value = -6249088
$ws.Range("N105").Value = -4984.2

$ws.Range("H107").Value = 1344.619
$ws.Range("I107").Value = 593.8182
$ws.Range("J107").Value = 2170.5
$ws.Range("K107").Value = 593.8182
$ws.Range("L107").Value = 2170.5
$ws.Range("M107").Value = 1326.1818
$ws.Range("N107").Value = -6010.5

$ws.Range("H113").Value = 1222.9166
$ws.Range("I113").Value = 1051.3334
$ws.Range("J113").Value = 1737.6666
$ws.Range("K113").Value = 1051.3334
$ws.Range("L113").Value = 1737.6666
$ws.Range("M113").Value = 1118.6666
$ws.Range("N113").Value = -6077.6666

$ws.Range("H126").Value = 33336580
$ws.Range("I126").Value = 2244.5715
$ws.Range("J126").Value = 62504124
$ws.Range("K126").Value = 6733.7145
$ws.Range("L126").Value = 187512372
$ws.Range("M126").Value = -4263.7145
$ws.Range("N126").Value = -187517312

$ws.Range("H132").Value = 2329
$ws.Range("I132").Value = 1572.12
$ws.Range("J132").Value = 4694.25
$ws.Range("K132").Value = 4716.36
$ws.Range("L132").Value = 14082.75
$ws.Range("M132").Value = -2186.36
$ws.Range("N132").Value = -19142.75

$ws.Range("H134").Value = 953.1429000000001
$ws.Range("I134").Value = 847.08
$ws.Range("J134").Value = 1218.3
$ws.Range("K134").Value = 2541.24
$ws.Range("L134").Value = 3654.9
$ws.Range("M134").Value = -6.240000000000236
$ws.Range("N134").Value = -8724.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1208.125
$ws.Range("I5").Value = 630.7778
$ws.Range("J5").Value = 1950.4286
$ws.Range("K5").Value = 1892.3334
$ws.Range("L5").Value = 5851.2858
$ws.Range("M5").Value = -1780.3334
$ws.Range("N5").Value = -6075.2858

$ws.Range("H12").Value = 29.266666
$ws.Range("I12").Value = 9.428572000000001
$ws.Range("J12").Value = 46.625
$ws.Range("K12").Value = 28.285716
$ws.Range("L12").Value = 139.875
$ws.Range("M12").Value = 144.714284
$ws.Range("N12").Value = -485.875

$ws.Range("H98").Value = 5000
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 5000
$ws.Range("K98").Value = 0
$ws.Range("L98").ClearContents()
$ws.Range("N98").Value = -17996

$ws.Range("H131").Value = 700.38
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 700.38
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2101.14
$ws.Range("N131").Value = -12181.14

$ws.Range("H135").Value = 1208.125
$ws.Range("I135").Value = 630.7778
$ws.Range("J135").Value = 1950.4286
$ws.Range("K135").Value = 5677.000199999999
$ws.Range("L135").Value = 17553.8574
$ws.Range("M135").Value = -3142.000199999999
$ws.Range("N135").Value = -22623.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18714.285
$ws.Range("I70").Value = 5500
$ws.Range("J70").Value = 36333.332
$ws.Range("K70").Value = 5500
$ws.Range("L70").Value = 36333.332
$ws.Range("M70").Value = -5230
$ws.Range("N70").Value = -36873.332

$ws.Range("H73").Value = 18714.285
$ws.Range("I73").Value = 5500
$ws.Range("J73").Value = 36333.332
$ws.Range("K73").Value = 5500
$ws.Range("L73").Value = 36333.332
$ws.Range("M73").Value = -4564
$ws.Range("N73").Value = -38205.332

$ws.Range("H102").Value = 14707312
$ws.Range("I102").Value = 16667819
$ws.Range("J102").Value = 3505.25
$ws.Range("K102").Value = 16667819
$ws.Range("L102").Value = 3505.25
$ws.Range("M102").Value = -16666197
$ws.Range("N102").Value = -6749.25

$ws.Range("H122").Value = 66669040
$ws.Range("I122").Value = 23810840
$ws.Range("J122").Value = 166671500
$ws.Range("K122").Value = 71432520
$ws.Range("L122").Value = 500014500
$ws.Range("M122").Value = -71430070
$ws.Range("N122").Value = -500019400

$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").ClearContents()

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").ClearContents()

$ws.Range("H132").Value = 14313.721
$ws.Range("I132").Value = 2602.359
$ws.Range("J132").Value = 128499.5
$ws.Range("K132").Value = 7807.076999999999
$ws.Range("L132").Value = 385498.5
$ws.Range("M132").Value = -5277.076999999999
$ws.Range("N132").Value = -390558.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2347.5
$ws.Range("I55").Value = 2347.5
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 2347.5
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()

$ws.Range("H61").Value = 3062.2104
$ws.Range("I61").Value = 1355.125
$ws.Range("J61").Value = 12166.667
$ws.Range("K61").Value = 1355.125
$ws.Range("L61").Value = 12166.667
$ws.Range("M61").Value = -1153.125
$ws.Range("N61").Value = -12570.667

$ws.Range("H113").Value = 3062.2104
$ws.Range("I113").Value = 1355.125
$ws.Range("J113").Value = 12166.667
$ws.Range("K113").Value = 1355.125
$ws.Range("L113").Value = 12166.667
$ws.Range("M113").Value = 814.875
$ws.Range("N113").Value = -16506.667

$ws.Range("H122").Value = 615535.5
$ws.Range("I122").Value = 982644.3
$ws.Range("J122").Value = 3687.5
$ws.Range("K122").Value = 2947932.9
$ws.Range("L122").Value = 11062.5
$ws.Range("M122").Value = -2945482.9
$ws.Range("N122").Value = -15962.5

$ws.Range("H132").Value = 484198.97
$ws.Range("I132").Value = 710182.1
$ws.Range("J132").Value = 3984.75
$ws.Range("K132").Value = 2130546.3
$ws.Range("L132").Value = 11954.25
$ws.Range("M132").Value = -2128016.3
$ws.Range("N132").Value = -17014.25

$ws.Range("H136").Value = 1554.6818
$ws.Range("I136").Value = 1424.3529
$ws.Range("J136").Value = 1997.8
$ws.Range("K136").Value = 4273.0587
$ws.Range("L136").Value = 5993.4
$ws.Range("M136").Value = -1723.0587
$ws.Range("N136").Value = -11093.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1537.9667
$ws.Range("I126").Value = 1182.1052
$ws.Range("J126").Value = 2152.6365
$ws.Range("K126").Value = 3546.3156
$ws.Range("L126").Value = 6457.9095
$ws.Range("M126").Value = -1076.3156
$ws.Range("N126").Value = -11397.9095

$ws.Range("H132").Value = 1950.4
$ws.Range("I132").Value = 1354.6154
$ws.Range("J132").Value = 3056.8572
$ws.Range("K132").Value = 4063.8462
$ws.Range("L132").Value = 9170.571599999999
$ws.Range("M132").Value = -1533.8462
$ws.Range("N132").Value = -14230.5716
